# Reorders the data rows (2-11) of Sheet1 per the "time bucket analysis"
# ordering added in this commit. The underlying data values themselves are
# unchanged - only the row each title/timestamp/historical-distance/
# time-bucket/uri combination lives on is different.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target row order (title, timestamp, historical distance, time bucket, uri)
$rows = @(
  @("Desert storm hits Taranaki", "2009-09-25T10:40:00UTC", 24, "day_2_to_30", "http://www.stuff.co.nz/national/2901968/Desert-storm-hits-Taranaki"),
  @("Are the dust storms radioactive? Australian scientists study Aussie dust from New Zealand", "2009-09-26T00:00:00UTC", 25, "day_2_to_30", "https://web.archive.org/web/20090928034025/http://www.news.com.au/story/0,27574,26127235-421,00.html"),
  @("Melbourne Airport warns of delays", "2009-09-23T08:12:00UTC", 22, "day_2_to_30", "http://www.abc.net.au/news/stories/2009/09/23/2693764.htm"),
  @("Red dust: face masks flying off the shelves", "2009-09-23T01:47:00UTC", 22, "day_2_to_30", "http://www.smh.com.au/environment/red-dust-face-masks-flying-off-the-shelves-20090923-g1jc.html"),
  @("Dust Storm Hits Brisbane, Queensland", "2009-09-23T00:00:00UTC", 22, "day_2_to_30", "https://web.archive.org/web/20090930184642/http://www.brisbanetimes.com.au/queensland/brisbane-awakes-from-haze-as-dust-settles-20090923-g0wm.html"),
  @("Red dust cloud covers Sydney", "2009-09-23T06:19:29UTC", 22, "day_2_to_30", "https://web.archive.org/web/20090923100531/http://news.bbc.co.uk/2/hi/asia-pacific/8270125.stm"),
  @("Red dust covers Sydney", "2009-09-23T06:19:29UTC", 22, "day_2_to_30", "http://news.bbc.co.uk/2/hi/asia-pacific/8270125.stm"),
  @("Canberra disappears in the dust", "2009-09-22T00:00:00UTC", 21, "day_2_to_30", "http://www.abc.net.au/news/photos/2009/09/22/2693220.htm"),
  @("Brisbane CBD dusted: Sydney storm heads north", "2009-09-23T00:00:00UTC", 22, "day_2_to_30", "http://www.theage.com.au/environment/brisbane-cbd-dusted-sydney-storm-heads-north-20090923-g0y1.html?autostart=1"),
  @("Sydney dust storm worst in 70 years, says weather bureau", "1-01-01T00:00:00UTC", "unknown", "unknown", "https://web.archive.org/web/20090927153112/http://www.theaustralian.news.com.au/story/0,25197,26113952-5006784,00.html")
)

# Write the reordered title / timestamp / historical distance / time bucket
# / uri columns for each row.
for ($i = 0; $i -lt $rows.Count; $i++) {
  $r = $i + 2
  $row = $rows[$i]

  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
}

# Drop the existing hyperlinks, then recreate them - in the SAME relative
# order the original links were defined in (news.bbc.co.uk, stuff.co.nz,
# archived-news.bbc.co.uk, abc stories, smh, abc photos, brisbanetimes,
# archived-news.com.au, theage, theaustralian) - so the underlying
# relationship ids line up with the original workbook, just retargeted at
# whichever row now holds that article.
$ws.Range("E2:E11").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Cells.Item(8, 5), "http://news.bbc.co.uk/2/hi/asia-pacific/8270125.stm")
$ws.Hyperlinks.Add($ws.Cells.Item(2, 5), "http://www.stuff.co.nz/national/2901968/Desert-storm-hits-Taranaki")
$ws.Hyperlinks.Add($ws.Cells.Item(7, 5), "https://web.archive.org/web/20090923100531/http://news.bbc.co.uk/2/hi/asia-pacific/8270125.stm")
$ws.Hyperlinks.Add($ws.Cells.Item(4, 5), "http://www.abc.net.au/news/stories/2009/09/23/2693764.htm")
$ws.Hyperlinks.Add($ws.Cells.Item(5, 5), "http://www.smh.com.au/environment/red-dust-face-masks-flying-off-the-shelves-20090923-g1jc.html")
$ws.Hyperlinks.Add($ws.Cells.Item(9, 5), "http://www.abc.net.au/news/photos/2009/09/22/2693220.htm")
$ws.Hyperlinks.Add($ws.Cells.Item(6, 5), "https://web.archive.org/web/20090930184642/http://www.brisbanetimes.com.au/queensland/brisbane-awakes-from-haze-as-dust-settles-20090923-g0wm.html")
$ws.Hyperlinks.Add($ws.Cells.Item(3, 5), "https://web.archive.org/web/20090928034025/http://www.news.com.au/story/0,27574,26127235-421,00.html")
$ws.Hyperlinks.Add($ws.Cells.Item(10, 5), "http://www.theage.com.au/environment/brisbane-cbd-dusted-sydney-storm-heads-north-20090923-g0y1.html?autostart=1")
$ws.Hyperlinks.Add($ws.Cells.Item(11, 5), "https://web.archive.org/web/20090927153112/http://www.theaustralian.news.com.au/story/0,25197,26113952-5006784,00.html")

# Restore the Hyperlink cell style (font/underline) on the uri column.
for ($r = 2; $r -le 11; $r++) {
  $ws.Cells.Item($r, 5).Style = "Hyperlink"
}
